# This script applies row-level data corrections to rows 32-40 of the
# "Artfynd" worksheet, matching the upstream automatic data refresh.
# For each affected cell we set the new value directly; string values that
# look like dates (columns Y/AA) are prefixed with a literal apostrophe so
# Excel stores them as text (matching the original inlineStr cells) rather
# than auto-converting them to date serial numbers. Cells that are removed
# in the target (K36, AC36) are cleared by assigning $null.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 112176074
$ws.Range("B32").Value = 98953
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 1365
$ws.Range("F32").Value = "Lappranunkel"
$ws.Range("G32").Value = "Coptidium lapponicum"
$ws.Range("H32").Value = "(L.) Tzvelev"
$ws.Range("Q32").Value = 602642
$ws.Range("R32").Value = 7030561

# Row 33
$ws.Range("B33").Value = 89834

# Row 34
$ws.Range("A34").Value = 112176108
$ws.Range("B34").Value = 89834
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 658
$ws.Range("F34").Value = "Rosenticka"
$ws.Range("G34").Value = "Rhodofomes roseus"
$ws.Range("H34").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q34").Value = 602831
$ws.Range("R34").Value = 7030665
$ws.Range("Y34").Value = "'2023-06-14"
$ws.Range("AA34").Value = "'2023-06-14"

# Row 36
$ws.Range("A36").Value = 112176087
$ws.Range("B36").Value = 89738
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 48
$ws.Range("F36").Value = "Lappticka"
$ws.Range("G36").Value = "Amylocystis lapponica"
$ws.Range("H36").Value = "(Romell) Singer"
$ws.Range("K36").Value = $null
$ws.Range("Q36").Value = 602806
$ws.Range("R36").Value = 7030689
$ws.Range("Y36").Value = "'2023-06-27"
$ws.Range("AA36").Value = "'2023-06-27"
$ws.Range("AC36").Value = $null
$ws.Range("AX36").Value = "Pekka Bader, Anna-Maria Eriksson"

# Row 37
$ws.Range("A37").Value = 112176095
$ws.Range("B37").Value = 89993
$ws.Range("E37").Value = 1209
$ws.Range("F37").Value = "Rynkskinn"
$ws.Range("G37").Value = "Phlebia centrifuga"
$ws.Range("H37").Value = "P.Karst."
$ws.Range("Q37").Value = 602796
$ws.Range("R37").Value = 7030566
$ws.Range("Y37").Value = "'2023-06-14"
$ws.Range("AA37").Value = "'2023-06-14"
$ws.Range("AX37").Value = "Pekka Bader"

# Row 38
$ws.Range("A38").Value = 112176093
$ws.Range("B38").Value = 6203
$ws.Range("E38").Value = 105336
$ws.Range("F38").Value = "Vanlig flatbagge"
$ws.Range("G38").Value = "Peltis ferruginea"
$ws.Range("H38").Value = "(Linnaeus, 1758)"
$ws.Range("Q38").Value = 602865
$ws.Range("R38").Value = 7030578
$ws.Range("Y38").Value = "'2023-06-22"
$ws.Range("AA38").Value = "'2023-06-22"
$ws.Range("AX38").Value = "Pekka Bader"

# Row 39
$ws.Range("A39").Value = 112176096
$ws.Range("B39").Value = 12450
$ws.Range("D39").Value = "EN"
$ws.Range("E39").Value = 101692
$ws.Range("F39").Value = "Större barkplattbagge"
$ws.Range("G39").Value = "Pytho kolwensis"
$ws.Range("H39").Value = "Sahlberg, 1833"
$ws.Range("K39").Value = "larv/nymf"
$ws.Range("Q39").Value = 602869
$ws.Range("R39").Value = 7030590
$ws.Range("AC39").Value = "larv 20-25 mm"

# Row 40
$ws.Range("A40").Value = 112176088
$ws.Range("B40").Value = 89573
$ws.Range("E40").Value = 5442
$ws.Range("F40").Value = "Tallticka"
$ws.Range("G40").Value = "Porodaedalea pini"
$ws.Range("H40").Value = "(Brot.) Murrill"
$ws.Range("Q40").Value = 602859
$ws.Range("R40").Value = 7030591
$ws.Range("Y40").Value = "'2023-06-27"
$ws.Range("AA40").Value = "'2023-06-27"
$ws.Range("AX40").Value = "Pekka Bader, Anna-Maria Eriksson"
